$d = $word.ActiveDocument

# 1. Highlight the "Tree-based Data Structures (7 days)" run in green.
$rng1 = $d.Content
[void]$rng1.Find.Execute("Tree-based Data Structures (7 days)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Font.HighlightColorIndex = 4

# 2. Highlight the " (06/07)" run in green as well, then insert a plain space
#    right after it.
$rng2 = $d.Content
[void]$rng2.Find.Execute(" (06/07)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Font.HighlightColorIndex = 4
$rng2.Collapse(0)
$rng2.InsertAfter(" ")
$rng2.Collapse(0)

# 3. Insert "(04/07)" after that space and colour only this new text red.
$rng3 = $d.Content
[void]$rng3.Find.Execute(" (06/07) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Collapse(0)
$rng3.InsertAfter("(04/07)")

$rng4 = $d.Content
[void]$rng4.Find.Execute("(04/07)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng4.Font.Color = 255
